# Updated cryptos list on Tue Dec  5 07:50:11 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matching source formatting),
# since the Price column stores values like "229.24" as plain text, not numbers.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '41.586.23'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '2.202.49'
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '229.24'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').Value = '0.617'
$ws.Range('E6').Value = '  -3.18%  '
$ws.Range('D7').Value = '59.78'
$ws.Range('E7').Value = '  -6.48%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.401'
$ws.Range('E9').Value = '  -2.10%  '
$ws.Range('D10').Value = '57.42'
$ws.Range('E10').Value = '  -3.36%  '
$ws.Range('D11').Value = '0.0888'
$ws.Range('E11').Value = '  -1.42%  '
$ws.Range('E12').Value = '  -1.15%  '
$ws.Range('D13').Value = '2.525.99'
$ws.Range('E13').Value = '  -2.33%  '
$ws.Range('E14').Value = '  -5.05%  '
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').Value = '5.64'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '0.793'
$ws.Range('E17').Value = '  -4.17%  '
$ws.Range('D18').Value = '2.255.31'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').Value = '41.437.80'
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').Value = '0.0₃0899'
$ws.Range('E20').Value = '  -2.36%  '
$ws.Range('D21').Value = '71.78'
$ws.Range('E21').Value = '  -2.52%  '
$ws.Range('D22').Value = '6.05'
$ws.Range('E22').Value = '  -2.15%  '
$ws.Range('D23').Value = '242.60'
$ws.Range('E23').Value = '  -3.36%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.36'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = '2.34'
$ws.Range('E26').Value = '  -1.86%  '
$ws.Range('D27').Value = '9.66'
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('D28').Value = '168.38'
$ws.Range('E28').Value = '  -2.72%  '
$ws.Range('E29').Value = '  -4.88%  '
$ws.Range('D30').Value = '19.70'
$ws.Range('E30').Value = '  -3.73%  '
$ws.Range('E31').Value = '  -2.57%  '
$ws.Range('E32').Value = '  -10.01%  '
$ws.Range('E33').Value = '  -2.96%  '
$ws.Range('D34').Value = '4.96'
$ws.Range('E34').Value = '  -1.88%  '
$ws.Range('D35').Value = '4.62'
$ws.Range('E35').Value = '  -2.35%  '
$ws.Range('D36').Value = '0.0645'
$ws.Range('E36').Value = '  +1.73%  '
$ws.Range('D37').Value = '6.44'
$ws.Range('E37').Value = '  -7.27%  '
$ws.Range('E38').Value = '  -3.81%  '
$ws.Range('D39').Value = '3.55'
$ws.Range('E39').Value = '  -7.61%  '
$ws.Range('D40').Value = '0.997'
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').Value = '0.000234'
$ws.Range('E41').Value = '  -13.50%  '
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('D43').Value = '8.53'
$ws.Range('E43').Value = '  -3.39%  '
$ws.Range('D44').Value = '0.0951'
$ws.Range('E44').Value = '  +0.80%  '
$ws.Range('E45').Value = '  -2.57%  '
$ws.Range('D46').Value = '97.04'
$ws.Range('E46').Value = '  -5.29%  '
$ws.Range('D47').Value = '1.463.27'
$ws.Range('E47').Value = '  -3.24%  '
$ws.Range('D48').Value = '4.32'
$ws.Range('E48').Value = '  -11.93%  '
$ws.Range('D49').Value = '16.35'
$ws.Range('E49').Value = '  -7.57%  '
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('D51').Value = '1.06'
$ws.Range('E51').Value = '  -4.95%  '
